# Fill in SqlServer row data and update the active selection, per commit
# "modified for server's data".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ID, ServerID, IP, Port, Pwd
# (written in this order so new shared-string entries land in the same
# table order as the target workbook: SqlServer_1, 127.0.0.1, 000107001)
$ws.Range("A2").Value = "SqlServer_1"
$ws.Range("C2").Value = "127.0.0.1"
$ws.Range("B2").Value = "000107001"
$ws.Range("D2").Value = 7001
$ws.Range("E2").Value = 123456

# Update selection to match the saved workbook state (active cell E4).
$ws.Range("E4").Select()
